# Apply "custom accuracy" rounding to row 5 data values (round to 2 decimals),
# shrink the width of columns G, K and Q by one unit (8 -> 7), and drop the
# last data row (row 6), shrinking the sheet from A1:AH6 to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: round the numeric readings to 2 decimal places -----------------
$ws.Range("B5").Value = 17.16
$ws.Range("C5").Value = 12.84
$ws.Range("D5").Value = 1.09
$ws.Range("E5").Value = 37.62
$ws.Range("F5").Value = 30.6
$ws.Range("G5").Value = 13.16
$ws.Range("H5").Value = 54.27
$ws.Range("I5").Value = 20.86
$ws.Range("J5").Value = 9.49
$ws.Range("K5").Value = 13.5
$ws.Range("L5").Value = 15.09
$ws.Range("M5").Value = 16.09
$ws.Range("N5").Value = 4.37
$ws.Range("O5").Value = 13.53
$ws.Range("P5").Value = 19.14
$ws.Range("Q5").Value = 11.52
$ws.Range("R5").Value = 0.33
$ws.Range("S5").Value = 0.68
$ws.Range("T5").Value = 198.7
$ws.Range("U5").Value = 37.89
$ws.Range("V5").Value = 12.49
$ws.Range("W5").Value = 25.35
$ws.Range("X5").Value = 13.47
$ws.Range("Y5").Value = 1.79
$ws.Range("Z5").Value = 26.7
$ws.Range("AA5").Value = 11.04
$ws.Range("AB5").Value = 9.82
$ws.Range("AC5").Value = 11.51
$ws.Range("AD5").Value = 15.92
$ws.Range("AE5").Value = 0.47
$ws.Range("AF5").Value = 49.35
$ws.Range("AG5").Value = 6.98
$ws.Range("AH5").Value = 15.62

# --- Remove row 6 entirely (last sample), sheet data shrinks to row 5 ------
$ws.Rows.Item(6).Delete()

# --- Narrow columns G (7), K (11) and Q (17) from width 8 to width 7 -------
# (ColumnWidth is expressed in character-width units; 6.14 lands in the
# quantization bucket that Excel persists as a stored width of exactly 7,
# matching the other width-7 columns already in the sheet.)
$ws.Columns.Item(7).ColumnWidth = 6.14
$ws.Columns.Item(11).ColumnWidth = 6.14
$ws.Columns.Item(17).ColumnWidth = 6.14
